$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "35.318.63"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "1.845.10"
$ws.Range("E3").Value = "  +1.86%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.23%  "
Set-TextValue "D5" "228.33"
$ws.Range("E5").Value = "  +1.20%  "
Set-TextValue "D6" "0.612"
$ws.Range("E6").Value = "  +2.06%  "
Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  -0.23%  "
Set-TextValue "D8" "42.99"
$ws.Range("E8").Value = "  +15.06%  "
Set-TextValue "D9" "0.306"
$ws.Range("E9").Value = "  +4.21%  "
Set-TextValue "D10" "0.0693"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").Value = "2.113.47"
$ws.Range("E12").Value = "  +1.91%  "
Set-TextValue "D13" "11.57"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "1.841.33"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E15").Value = "  +7.08%  "
Set-TextValue "D16" "0.660"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "35.158.24"
$ws.Range("E17").Value = "  +1.76%  "
Set-TextValue "D18" "69.73"
$ws.Range("E18").Value = "  +1.51%  "
Set-TextValue "D19" "246.70"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  +2.20%  "
Set-TextValue "D21" "12.09"
$ws.Range("E21").Value = "  +7.52%  "
$ws.Range("E22").Value = "  +13.99%  "
$ws.Range("E23").Value = "  -0.17%  "
Set-TextValue "D24" "2.19"
$ws.Range("E24").Value = "  -1.00%  "
Set-TextValue "D25" "171.54"
$ws.Range("E25").Value = "  -0.33%  "
Set-TextValue "D26" "7.94"
$ws.Range("E26").Value = "  +0.66%  "
Set-TextValue "D27" "17.90"
$ws.Range("E27").Value = "  +3.51%  "
Set-TextValue "D28" "0.123"
$ws.Range("E28").Value = "  +1.01%  "
Set-TextValue "D29" "0.997"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +8.82%  "
Set-TextValue "D31" "3.96"
$ws.Range("E31").Value = "  +3.49%  "
Set-TextValue "D32" "4.05"
$ws.Range("E32").Value = "  +2.76%  "
Set-TextValue "D33" "0.0537"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("E34").Value = "  +4.41%  "
Set-TextValue "D35" "0.673"
$ws.Range("E35").Value = "  +2.90%  "
Set-TextValue "D36" "90.33"
$ws.Range("E36").Value = "  +11.28%  "
Set-TextValue "D37" "1.08"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "1.340.44"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  +9.10%  "
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("E42").Value = "  +8.24%  "
$ws.Range("E43").Value = "  +6.54%  "
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  +3.23%  "
Set-TextValue "D47" "6.06"
$ws.Range("E47").Value = "  +4.13%  "
$ws.Range("D48").Value = "2.012.39"
$ws.Range("E48").Value = "  +1.95%  "
Set-TextValue "D49" "104.47"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D51" "49.39"
$ws.Range("E51").Value = "  +2.06%  "
